# Auto-generated script to apply scheduled-runner price/profit updates
# to the Cactuar_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 47621520
$ws.Range("J18").Value = 2800
$ws.Range("L18").Value = 2800
$ws.Range("N18").Value = -3368

$ws.Range("H43").Value = 5133038
$ws.Range("I43").Value = 15384615
$ws.Range("K43").Value = 15384615
$ws.Range("M43").Value = -15384546

$ws.Range("H100").Value = 2883.5
$ws.Range("I100").Value = 2973.5
$ws.Range("J100").Value = 2823.5
$ws.Range("K100").Value = 2973.5
$ws.Range("L100").Value = 2823.5
$ws.Range("M100").Value = -2432.5
$ws.Range("N100").Value = -3905.5

$ws.Range("H106").Value = 47621290
$ws.Range("I106").Value = 55557836
$ws.Range("K106").Value = 55557836
$ws.Range("M106").Value = -55557205

$ws.Range("H116").Value = 50699616
$ws.Range("I116").Value = 101395890
$ws.Range("J116").Value = 3347.3333
$ws.Range("K116").Value = 101395890
$ws.Range("L116").Value = 3347.3333
$ws.Range("M116").Value = -101392448
$ws.Range("N116").Value = -10231.3333

$ws.Range("H129").Value = 2514.75
$ws.Range("J129").Value = 4060
$ws.Range("L129").Value = 12180
$ws.Range("N129").Value = -22180

$ws.Range("H132").Value = 5942.689
$ws.Range("I132").Value = 1220.9445
$ws.Range("K132").Value = 3662.8335
$ws.Range("M132").Value = -1132.8335

$ws.Range("H141").Value = 12375.207
$ws.Range("J141").Value = 14066.333
$ws.Range("L141").Value = 42198.999
$ws.Range("N141").Value = -52558.999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 805.1667
$ws.Range("I4").Value = 226.33333
$ws.Range("J4").Value = 1384
$ws.Range("K4").Value = 226.33333
$ws.Range("L4").Value = 1384
$ws.Range("M4").Value = -110.33333
$ws.Range("N4").Value = -1616

$ws.Range("H32").Value = 25422.29
$ws.Range("I32").Value = 25439.062
$ws.Range("K32").Value = 25439.062
$ws.Range("M32").Value = -25152.062

$ws.Range("H45").Value = 3034.75
$ws.Range("I45").Value = 1999.6
$ws.Range("K45").Value = 1999.6
$ws.Range("M45").Value = -1622.6

$ws.Range("H61").Value = 4342.196
$ws.Range("I61").Value = 3650.9412
$ws.Range("K61").Value = 3650.9412
$ws.Range("M61").Value = -3438.9412

$ws.Range("H122").Value = 7620.6665
$ws.Range("I122").Value = 6981.8
$ws.Range("J122").Value = 8201.454
$ws.Range("K122").Value = 20945.4
$ws.Range("L122").Value = 24604.362
$ws.Range("M122").Value = -18495.4
$ws.Range("N122").Value = -29504.362

$ws.Range("H136").Value = 4342.196
$ws.Range("I136").Value = 3650.9412
$ws.Range("K136").Value = 10952.8236
$ws.Range("M136").Value = -8402.8236

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3127.85
$ws.Range("I20").Value = 3672.111
$ws.Range("J20").Value = 2682.5454
$ws.Range("K20").Value = 3672.111
$ws.Range("L20").Value = 2682.5454
$ws.Range("M20").Value = -3425.111
$ws.Range("N20").Value = -3176.5454

$ws.Range("H22").Value = 486.44446
$ws.Range("I22").Value = 297
$ws.Range("J22").Value = 2002
$ws.Range("K22").Value = 297
$ws.Range("L22").Value = 2002
$ws.Range("M22").Value = -124
$ws.Range("N22").Value = -2348

$ws.Range("H105").Value = 4999.6665
$ws.Range("I105").Value = 4000
$ws.Range("K105").Value = 4000
$ws.Range("M105").Value = -2253

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11496274
$ws.Range("I31").Value = 13159466
$ws.Range("K31").Value = 13159466
$ws.Range("M31").Value = -13159171

$ws.Range("H34").Value = 11496274
$ws.Range("I34").Value = 13159466
$ws.Range("K34").Value = 13159466
$ws.Range("M34").Value = -13159264

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()

$ws.Range("H130").Value = 85000
$ws.Range("J130").Value = 85000
$ws.Range("L130").Value = 85000
$ws.Range("N130").Value = -95040

$ws.Range("H132").Value = 74082320
$ws.Range("I132").Value = 102566216
$ws.Range("J132").Value = 24189.8
$ws.Range("K132").Value = 307698648
$ws.Range("L132").Value = 72569.39999999999
$ws.Range("M132").Value = -307696118
$ws.Range("N132").Value = -77629.39999999999

$ws.Range("H134").Value = 3399.4
$ws.Range("I134").Value = 2699.5
$ws.Range("K134").Value = 8098.5
$ws.Range("M134").Value = -5563.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1924.5
$ws.Range("I22").Value = 2466.6667
$ws.Range("J22").Value = 1599.2
$ws.Range("K22").Value = 7400.000100000001
$ws.Range("L22").Value = 4797.6
$ws.Range("M22").Value = -7231.000100000001
$ws.Range("N22").Value = -5135.6

$ws.Range("H27").Value = 1924.5
$ws.Range("I27").Value = 2466.6667
$ws.Range("J27").Value = 1599.2
$ws.Range("K27").Value = 7400.000100000001
$ws.Range("L27").Value = 4797.6
$ws.Range("M27").Value = -7298.000100000001
$ws.Range("N27").Value = -5001.6

$ws.Range("H118").Value = 5601.231
$ws.Range("I118").Value = 5488
$ws.Range("J118").Value = 5733.3335
$ws.Range("K118").Value = 16464
$ws.Range("L118").Value = 17200.0005
$ws.Range("M118").Value = -15221
$ws.Range("N118").Value = -19686.0005

$ws.Range("H131").Value = 24981348
$ws.Range("J131").Value = 26125356
$ws.Range("L131").Value = 78376068
$ws.Range("N131").Value = -78386148

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2074051.1
$ws.Range("I70").Value = 3795494.5
$ws.Range("K70").Value = 3795494.5
$ws.Range("M70").Value = -3795224.5

$ws.Range("H73").Value = 2074051.1
$ws.Range("I73").Value = 3795494.5
$ws.Range("K73").Value = 3795494.5
$ws.Range("M73").Value = -3794558.5

$ws.Range("H88").Value = 120000
$ws.Range("J88").Value = 120000
$ws.Range("L88").Value = 120000
$ws.Range("N88").Value = -120902

$ws.Range("H91").Value = 120000
$ws.Range("J91").Value = 120000
$ws.Range("L91").Value = 120000
$ws.Range("N91").Value = -123120

$ws.Range("H122").Value = 421400.66
$ws.Range("I122").Value = 835801.75
$ws.Range("J122").Value = 6999.5835
$ws.Range("K122").Value = 2507405.25
$ws.Range("L122").Value = 20998.7505
$ws.Range("M122").Value = -2504955.25
$ws.Range("N122").Value = -25898.7505

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1057.5264
$ws.Range("I16").Value = 618.26666
$ws.Range("J16").Value = 2704.75
$ws.Range("K16").Value = 618.26666
$ws.Range("L16").Value = 2704.75
$ws.Range("M16").Value = -448.26666
$ws.Range("N16").Value = -3044.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 84300
$ws.Range("J123").Value = 84300
$ws.Range("L123").Value = 84300
$ws.Range("N123").Value = -94100

$ws.Range("H132").Value = 4821.636
$ws.Range("I132").Value = 4749.143
$ws.Range("J132").Value = 4948.5
$ws.Range("K132").Value = 14247.429
$ws.Range("L132").Value = 14845.5
$ws.Range("M132").Value = -11717.429
$ws.Range("N132").Value = -19905.5

$ws.Range("H136").Value = 2202.0605
$ws.Range("I136").Value = 2124.8064
$ws.Range("K136").Value = 6374.4192
$ws.Range("M136").Value = -3824.4192
